$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition) - update "want to go" counts (column F) only
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 1563
$wsExpo.Range("F5").Value  = 8745
$wsExpo.Range("F8").Value  = 1241
$wsExpo.Range("F10").Value = 241
$wsExpo.Range("F13").Value = 107
$wsExpo.Range("F14").Value = 279
$wsExpo.Range("F15").Value = 7
$wsExpo.Range("F17").Value = 1405
$wsExpo.Range("F18").Value = 1306
$wsExpo.Range("F21").Value = 1327
$wsExpo.Range("F25").Value = 66
$wsExpo.Range("F27").Value = 276
$wsExpo.Range("F31").Value = 202
$wsExpo.Range("F35").Value = 597
$wsExpo.Range("F39").Value = 143
$wsExpo.Range("F41").Value = 1218

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance) - update "want to go" counts (column F) only
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F7").Value  = 42
$wsShow.Range("F8").Value  = 9
$wsShow.Range("F24").Value = 919
$wsShow.Range("F27").Value = 167

# ---------------------------------------------------------------------------
# Sheet "本地生活" (Local Life) - the "Paradox Live" event (row 4) was removed
# entirely, shifting all following rows up by one. Some of those shifted
# rows also received updated "want to go" counts.
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Rows.Item(4).Delete() | Out-Null

# Fix up the sequential index column (A) for the rows that shifted up.
$wsLocal.Range("A4").Value = 3
$wsLocal.Range("A5").Value = 4
$wsLocal.Range("A6").Value = 5
$wsLocal.Range("A7").Value = 6
$wsLocal.Range("A8").Value = 7
$wsLocal.Range("A9").Value = 8

# Updated "want to go" counts among the shifted rows.
$wsLocal.Range("F5").Value = 730
$wsLocal.Range("F7").Value = 129
$wsLocal.Range("F8").Value = 1957
$wsLocal.Range("F9").Value = 2969

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - update several "want to go" counts (column F)
# and then insert a new row for "上海·反派角色only展" at row 29, shifting all
# later rows down by one.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 730
$wsAll.Range("F7").Value  = 8745
$wsAll.Range("F9").Value  = 129
$wsAll.Range("F10").Value = 9
$wsAll.Range("F11").Value = 1957
$wsAll.Range("F12").Value = 2969
$wsAll.Range("F17").Value = 1241
$wsAll.Range("F21").Value = 107
$wsAll.Range("F22").Value = 279
$wsAll.Range("F23").Value = 1405
$wsAll.Range("F24").Value = 1306
$wsAll.Range("F25").Value = 1327
$wsAll.Range("F28").Value = 276

# Insert the new row for "上海·反派角色only展"
$wsAll.Rows.Item(29).Insert() | Out-Null

# Fix up the sequential index column (A) for the rows that shifted down.
for ($r = 30; $r -le 50; $r++) {
    $wsAll.Cells.Item($r, 1).Value = $r - 1
}

$wsAll.Range("A29").Value = 28
$wsAll.Range("B29").Value = "2024-06-09"
$wsAll.Range("C29").Value = "上海·反派角色only展"
$wsAll.Range("D29").Value = "长江路 258号中成智谷创意产业园区 成美术馆"
$wsAll.Range("E29").Value = "2024.06.09 09:30-06.09 17:30"
$wsAll.Range("F29").Value = 276
$wsAll.Range("G29").Value = 89.1
$wsAll.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=85071"
$wsAll.Range("I29").Value = "//i1.hdslb.com/bfs/openplatform/202404/gOvw8Iip1714470948525.jpeg"
